$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add two new "Stretch Goals" entries in G/H columns for rows 12 and 13
$ws.Range("G12").Value = "Player attacked animation"
$ws.Range("H12").Value = 3
$ws.Range("G13").Value = "Ghost attacked animation"
$ws.Range("H13").Value = 3

# Mark rows 28 and 29 as "WIP" in the Finished? column (E)
$ws.Range("E28").Value = "WIP"
$ws.Range("E29").Value = "WIP"

# Change row 32 Finished? from "Yes" to "WIP"
$ws.Range("E32").Value = "WIP"

# Update the view so it matches the saved window/selection state
# (scroll so row 22 is the top visible row, column A stays leftmost)
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E32").Select() | Out-Null
